# Update the multiplication-table answers to a new randomized problem set.
$d = $word.ActiveDocument

$d.Content.Find.Execute("552×4=2208", $true, $false, $false, $false, $false, $true, 1, $false, "924×8=7392", 2)
$d.Content.Find.Execute("282×7=1974", $true, $false, $false, $false, $false, $true, 1, $false, "157×3=471", 2)
$d.Content.Find.Execute("200×8=1600", $true, $false, $false, $false, $false, $true, 1, $false, "136×6=816", 2)
$d.Content.Find.Execute("148×7=1036", $true, $false, $false, $false, $false, $true, 1, $false, "708×7=4956", 2)
$d.Content.Find.Execute("920×2=1840", $true, $false, $false, $false, $false, $true, 1, $false, "589×3=1767", 2)
$d.Content.Find.Execute("993×3=2979", $true, $false, $false, $false, $false, $true, 1, $false, "722×3=2166", 2)
$d.Content.Find.Execute("271×5=1355", $true, $false, $false, $false, $false, $true, 1, $false, "776×7=5432", 2)
$d.Content.Find.Execute("991×8=7928", $true, $false, $false, $false, $false, $true, 1, $false, "767×2=1534", 2)
$d.Content.Find.Execute("467×7=3269", $true, $false, $false, $false, $false, $true, 1, $false, "308×7=2156", 2)
$d.Content.Find.Execute("591×2=1182", $true, $false, $false, $false, $false, $true, 1, $false, "372×2=744", 2)
$d.Content.Find.Execute("290×3=870", $true, $false, $false, $false, $false, $true, 1, $false, "161×8=1288", 2)
$d.Content.Find.Execute("505×9=4545", $true, $false, $false, $false, $false, $true, 1, $false, "803×6=4818", 2)
$d.Content.Find.Execute("939×3=2817", $true, $false, $false, $false, $false, $true, 1, $false, "950×4=3800", 2)
$d.Content.Find.Execute("300×2=600", $true, $false, $false, $false, $false, $true, 1, $false, "119×5=595", 2)
$d.Content.Find.Execute("581×7=4067", $true, $false, $false, $false, $false, $true, 1, $false, "484×5=2420", 2)
$d.Content.Find.Execute("623×2=1246", $true, $false, $false, $false, $false, $true, 1, $false, "949×4=3796", 2)
$d.Content.Find.Execute("569×8=4552", $true, $false, $false, $false, $false, $true, 1, $false, "211×4=844", 2)
$d.Content.Find.Execute("831×5=4155", $true, $false, $false, $false, $false, $true, 1, $false, "551×8=4408", 2)
$d.Content.Find.Execute("924×3=2772", $true, $false, $false, $false, $false, $true, 1, $false, "791×9=7119", 2)
$d.Content.Find.Execute("280×4=1120", $true, $false, $false, $false, $false, $true, 1, $false, "318×5=1590", 2)
$d.Content.Find.Execute("997×7=6979", $true, $false, $false, $false, $false, $true, 1, $false, "168×7=1176", 2)
$d.Content.Find.Execute("764×4=3056", $true, $false, $false, $false, $false, $true, 1, $false, "722×2=1444", 2)
$d.Content.Find.Execute("660×7=4620", $true, $false, $false, $false, $false, $true, 1, $false, "763×4=3052", 2)
$d.Content.Find.Execute("812×8=6496", $true, $false, $false, $false, $false, $true, 1, $false, "779×9=7011", 2)
$d.Content.Find.Execute("827×7=5789", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=3510", 2)
